# Auto-generated edit script: update crypto price/volume table (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.996.64"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.676.59"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.11"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "1.714.14"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.83"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.996.91"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.16"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "237.19"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.07"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.22"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.14"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "1.479.33"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.905"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.48"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.48"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.76"
$ws.Range("E51").Value = "  +1.20%  "
